# UndoRedoNewCommand1StackDiagram.pptx update:
# Rename the "prevTaskBook" identifier to "prevOrganizer" in every
# command-diagram table on the slide (PrioriTask was renamed from
# Addressbook/TaskBook to Organizer during the documentation pass
# described by the commit message).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable) {
        $tbl = $shape.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cell = $tbl.Cell($r, $c)
                $tr = $cell.Shape.TextFrame.TextRange
                $paras = $tr.Paragraphs()
                for ($pi = 1; $pi -le $paras.Count; $pi++) {
                    $para = $paras.Item($pi)
                    if ($para.Text -like "*prevTaskBook*") {
                        $para.Text = $para.Text.Replace("prevTaskBook", "prevOrganizer")
                    }
                }
            }
        }
    }
}
